$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 320 (pushes existing rows 320-396 down to 321-397)
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new data record
$ws.Range("A320").Value = 3
$ws.Range("B320").Value = "Femacal de La Calera"
$ws.Range("C320").Value = "Coquimbo"
$ws.Range("D320").Value = 44782
$ws.Range("E320").Value = 5
$ws.Range("F320").Value = 100112043
$ws.Range("G320").Value = "Pepino ensalada"
$ws.Range("H320").Value = "Sin especificar"
$ws.Range("I320").Value = "Primera"
$ws.Range("J320").Value = 90
$ws.Range("K320").Value = 18000
$ws.Range("L320").Value = 18500
$ws.Range("M320").Value = 18250
$ws.Range("N320").Value = "$/caja 70 unidades"
$ws.Range("O320").Value = "Región de Arica y Parinacota"
$ws.Range("P320").Value = 261
$ws.Range("Q320").Value = 70
$ws.Range("R320").Value = "Hortaliza"
